# Generate Report for Handback
# Updates the zh-cn and de-de sheets with the handback result of the
# 9c58bbe0-2c9c-499e-af8b-e0a1eefd73ba.md file: a "Latest Target File",
# "Latest Handback File", "Latest Handback DateTime" and an "Error Detail"
# describing that the handed-back file is not the latest version.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c2c1ce77e5151593b7a393f9291131b8ac78507/e2e/9c58bbe0-2c9c-499e-af8b-e0a1eefd73ba.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9666038e59b4898f5a1dc6b29c7fdb8ffe8c9c48/e2e/9c58bbe0-2c9c-499e-af8b-e0a1eefd73ba.md."
$mdDisplay = "9c58bbe0-2c9c-499e-af8b-e0a1eefd73ba.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c2c1ce77e5151593b7a393f9291131b8ac78507/e2e/9c58bbe0-2c9c-499e-af8b-e0a1eefd73ba.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen the "Latest Target File", "Latest Handback File" and
# "Error Detail" columns so the new long values are readable.
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# I6 - Latest Target File: add hyperlink to the handback markdown file
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(6, 9), $mdUrl, "", "", $mdDisplay)
$wsZh.Cells.Item(6, 9).Font.Underline = $true
$wsZh.Cells.Item(6, 9).Font.Color = 15570276

# J6 - Latest Handback File
$wsZh.Cells.Item(6, 10).Value = "9c58bbe0-2c9c-499e-af8b-e0a1eefd73ba.2eed646854b8f2afc0b0d2adefe2019387428899.zh-cn.xlf"

# K6 - Latest Handback DateTime
$wsZh.Cells.Item(6, 11).Value = "2016-10-13 12:59:11"

# P6 - Error Detail
$wsZh.Cells.Item(6, 16).Value = $errorDetail

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
$wsDe.Columns.Item(16).ColumnWidth = 39.17

# I6 - Latest Target File: add hyperlink to the handback markdown file
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(6, 9), $mdUrl, "", "", $mdDisplay)
$wsDe.Cells.Item(6, 9).Font.Underline = $true
$wsDe.Cells.Item(6, 9).Font.Color = 15570276

# J6 - Latest Handback File
$wsDe.Cells.Item(6, 10).Value = "9c58bbe0-2c9c-499e-af8b-e0a1eefd73ba.2eed646854b8f2afc0b0d2adefe2019387428899.de-de.xlf"

# K6 - Latest Handback DateTime
$wsDe.Cells.Item(6, 11).Value = "2016-10-13 12:59:26"

# P6 - Error Detail
$wsDe.Cells.Item(6, 16).Value = $errorDetail
